$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so price strings containing dots
# (e.g. "29.211.01") are not misread as numbers by the Excel parser.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '29.211.01'
$ws.Range('E2').Value = '  -1.09%  '
$ws.Range('D3').Value = '1.866.69'
$ws.Range('E3').Value = '  -0.80%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '0.7108'
$ws.Range('E5').Value = '  -0.63%  '
$ws.Range('D6').Value = '241.49'
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('D8').Value = '0.3118'
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').Value = '0.07666'
$ws.Range('E9').Value = '  -3.68%  '
$ws.Range('E10').Value = '  -2.19%  '
$ws.Range('D11').Value = '0.08369'
$ws.Range('E11').Value = '  +1.09%  '
$ws.Range('D12').Value = '1.873.40'
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('D13').Value = '5.225'
$ws.Range('D14').Value = '0.7116'
$ws.Range('E14').Value = '  -2.48%  '
$ws.Range('D15').Value = '91.34'
$ws.Range('E15').Value = '  +0.10%  '
$ws.Range('D16').Value = '29.237.81'
$ws.Range('E16').Value = '  -1.01%  '
$ws.Range('D17').Value = '5.948'
$ws.Range('E17').Value = '  -0.01%  '
$ws.Range('D18').Value = '243.65'
$ws.Range('E18').Value = '  -1.19%  '
$ws.Range('D19').Value = '0.000007825'
$ws.Range('E19').Value = '  -0.86%  '
$ws.Range('D20').Value = '2.113.67'
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('E21').Value = '  -2.06%  '
$ws.Range('D22').Value = '0.9993'
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('D23').Value = '7.863'
$ws.Range('E23').Value = '  -1.36%  '
$ws.Range('D25').Value = '0.1604'
$ws.Range('E25').Value = '  -0.87%  '
$ws.Range('D26').Value = '163.06'
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').Value = '8.960'
$ws.Range('E27').Value = '  -1.24%  '
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('E30').Value = '  -3.62%  '
$ws.Range('D31').Value = '4.402'
$ws.Range('E31').Value = '  +0.13%  '
$ws.Range('D32').Value = '4.248'
$ws.Range('E32').Value = '  +3.37%  '
$ws.Range('D33').Value = '0.05152'
$ws.Range('E33').Value = '  -2.46%  '
$ws.Range('D34').Value = '0.8004'
$ws.Range('E34').Value = '  +9.93%  '
$ws.Range('E35').Value = '  -2.66%  '
$ws.Range('E36').Value = '  -2.84%  '
$ws.Range('E37').Value = '  +0.20%  '
$ws.Range('D38').Value = '0.01855'
$ws.Range('E38').Value = '  -0.83%  '
$ws.Range('D39').Value = '2.709'
$ws.Range('E39').Value = '  -0.25%  '
$ws.Range('D40').Value = '1.159.58'
$ws.Range('E40').Value = '  -5.97%  '
$ws.Range('D41').Value = '6.312'
$ws.Range('E41').Value = '  +1.56%  '
$ws.Range('D42').Value = '0.8977'
$ws.Range('E42').Value = '  -1.50%  '
$ws.Range('D43').Value = '73.18'
$ws.Range('E43').Value = '  -1.09%  '
$ws.Range('D44').Value = '0.9999'
$ws.Range('E44').Value = '  -0.12%  '
$ws.Range('D45').Value = '103.21'
$ws.Range('E45').Value = '  +0.89%  '
$ws.Range('D46').Value = '2.011.57'
$ws.Range('E46').Value = '  -0.73%  '
$ws.Range('D47').Value = '0.5180'
$ws.Range('E47').Value = '  -2.07%  '
$ws.Range('E48').Value = '  -1.14%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '9.342'
$ws.Range('E49').Value = '  +0.18%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.00000000120'
$ws.Range('E50').Value = '  -0.54%  '
$ws.Range('D51').Value = '0.4297'
$ws.Range('E51').Value = '  -0.84%  '

# Restore default (General) formatting/style on column D so only
# the cell values changed, matching the original formatting.
$ws.Range('D2:D51').ClearFormats()
